# Fruta / hortaliza, semanal
# The weekly data refresh reshuffles the Fecha (D), Volumen (M),
# Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P)
# and Precio $/Kg (S) values across the existing data rows (2-9).
# All other columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move, keyed by row.
$cols = @("D", "M", "N", "O", "P", "S")
$before = @{}
for ($r = 2; $r -le 9; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: new row (key) receives the captured values that used to live
# in the row given as the value.
$mapping = @{
    2 = 7
    3 = 4
    4 = 2
    5 = 6
    6 = 8
    7 = 9
    8 = 5
    9 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $before[$srcRow][$c]
    }
}
